# fix number to string
# Cell B3 on the "configArray" sheet holds the declared type for the "age"
# column; it was "int" but should be "string".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configArray")

# Update the cell value (was shared-string "int" -> now "string")
$ws.Range("B3").Value = "string"

# Move the selection/active cell on this sheet, and make it the active tab.
$ws.Range("C13").Select()
$ws.Activate()
